$d = $word.ActiveDocument

$replacements = @(
    ,@("2023-10-23 Monday", "2023-10-24 Tuesday")
    ,@("42-23=", "70-4=")
    ,@("36+30=", "89-57=")
    ,@("91-20=", "92+7=")
    ,@("68-58=", "3+74=")
    ,@("70-64=", "15+72=")
    ,@("45+11=", "65+15=")
    ,@("76+9=", "57-25=")
    ,@("92-47=", "15+21=")
    ,@("21+36=", "15+28=")
    ,@("48+50=", "57+24=")
    ,@("98-22=", "0+74=")
    ,@("89-46=", "75-64=")
    ,@("50-25=", "95-82=")
    ,@("34+64=", "59-23=")
    ,@("8+34=", "43+4=")
    ,@("63+6=", "99-53=")
    ,@("25+68=", "46-4=")
    ,@("52-15=", "34+18=")
    ,@("14+64=", "26+44=")
    ,@("7+65=", "44+13=")
    ,@("12+17=", "67-31=")
    ,@("52+19=", "16+13=")
    ,@("36-5=", "96-9=")
    ,@("14+80=", "94-57=")
    ,@("96-68=", "44+5=")
    ,@("20-12=", "94-57=")
    ,@("81+11=", "56-41=")
    ,@("87-78=", "85-77=")
    ,@("33-26=", "5+35=")
    ,@("59+0=", "41+14=")
    ,@("89+2=", "85-63=")
    ,@("50+12=", "56-42=")
    ,@("38+29=", "71-64=")
    ,@("90-30=", "40+44=")
    ,@("92-91=", "24+38=")
    ,@("5+15=", "15+22=")
    ,@("3+11=", "33+19=")
    ,@("13+19=", "23+21=")
    ,@("58+0=", "6+46=")
    ,@("85-55=", "74-60=")
    ,@("44+32=", "48+21=")
    ,@("49+28=", "47-35=")
    ,@("62+10=", "30-19=")
    ,@("22+44=", "83+5=")
    ,@("29+23=", "21+71=")
    ,@("88-50=", "42-24=")
    ,@("70-27=", "43+19=")
    ,@("24+14=", "16+41=")
    ,@("43-17=", "79-39=")
    ,@("78-75=", "94-37=")
    ,@("11+59=", "59-17=")
    ,@("5+85=", "29+8=")
    ,@("17+3=", "71+23=")
    ,@("91-62=", "17+70=")
    ,@("83+14=", "71+23=")
    ,@("93-16=", "50-15=")
    ,@("48-39=", "89-36=")
    ,@("89-81=", "80-31=")
    ,@("66+27=", "95-93=")
    ,@("31+45=", "92-55=")
    ,@("70-32=", "8+50=")
    ,@("9+35=", "60-51=")
    ,@("66+1=", "32-10=")
    ,@("82-47=", "5+19=")
    ,@("85-61=", "76+23=")
    ,@("89-54=", "76-47=")
    ,@("75-75=", "6+90=")
    ,@("91-79=", "4+53=")
    ,@("1+17=", "68+31=")
    ,@("20+43=", "89-87=")
    ,@("17+10=", "76-67=")
    ,@("60+16=", "32+26=")
    ,@("49-47=", "48-6=")
    ,@("86-39=", "68-23=")
    ,@("83-65=", "35+26=")
    ,@("7+46=", "1+54=")
    ,@("66+10=", "94-7=")
    ,@("64+1=", "18+14=")
    ,@("44-3=", "53+33=")
    ,@("71-29=", "55-43=")
    ,@("17+29=", "33-8=")
    ,@("28+53=", "66+16=")
    ,@("93-69=", "80-14=")
    ,@("85-59=", "6+23=")
    ,@("73-68=", "40-39=")
    ,@("57+12=", "76-47=")
    ,@("25-12=", "39-31=")
    ,@("83-82=", "70+16=")
    ,@("71-48=", "7+54=")
    ,@("83-40=", "11-5=")
    ,@("85-71=", "95-29=")
    ,@("58-11=", "70-38=")
    ,@("0+34=", "58+13=")
    ,@("85-50=", "75-29=")
    ,@("93-31=", "72-69=")
    ,@("14+9=", "42+44=")
    ,@("19+17=", "48-11=")
    ,@("32+63=", "2+19=")
    ,@("14+49=", "85-34=")
    ,@("39+54=", "32+39=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Done applying replacements"
